$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 718.8570999999999
$ws.Range("J2").Value = 1053
$ws.Range("L2").Value = 1053
$ws.Range("N2").Value = -1279

$ws.Range("H52").Value = 600
$ws.Range("J52").Value = 600
$ws.Range("L52").Value = 1800
$ws.Range("N52").Value = -2120

$ws.Range("H105").Value = 47500
$ws.Range("J105").Value = 25000
$ws.Range("L105").Value = 25000
$ws.Range("N105").Value = -31988

$ws.Range("H106").Value = 18705.428
$ws.Range("I106").Value = 20373.143
$ws.Range("J106").Value = 13702.286
$ws.Range("K106").Value = 20373.143
$ws.Range("L106").Value = 13702.286
$ws.Range("M106").Value = -19742.143
$ws.Range("N106").Value = -14964.286

$ws.Range("H113").Value = 3499
$ws.Range("I113").Value = 3500
$ws.Range("K113").Value = 3500
$ws.Range("M113").Value = -246

$ws.Range("H131").Value = 750
$ws.Range("I131").Value = 750
$ws.Range("K131").Value = 2250
$ws.Range("M131").Value = 2790

$ws.Range("H137").Value = 1438.0714
$ws.Range("I137").Value = 1316.3334
$ws.Range("J137").Value = 1657.2
$ws.Range("K137").Value = 3949.0002
$ws.Range("L137").Value = 4971.6
$ws.Range("M137").Value = -1399.0002
$ws.Range("N137").Value = -10071.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 945
$ws.Range("I74").Value = 755.5454999999999
$ws.Range("K74").Value = 755.5454999999999
$ws.Range("M74").Value = 118.4545000000001

$ws.Range("H77").Value = 945
$ws.Range("I77").Value = 755.5454999999999
$ws.Range("K77").Value = 3777.7275
$ws.Range("M77").Value = 590.2725

$ws.Range("H122").Value = 486968.2
$ws.Range("I122").Value = 594494.9
$ws.Range("K122").Value = 1783484.7
$ws.Range("M122").Value = -1781034.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2414.0557
$ws.Range("J86").Value = 2620
$ws.Range("L86").Value = 2620
$ws.Range("N86").Value = -4866

$ws.Range("H89").Value = 2414.0557
$ws.Range("J89").Value = 2620
$ws.Range("L89").Value = 13100
$ws.Range("N89").Value = -24332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 4600
$ws.Range("I48").Value = 4600
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 4600
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -4124
$ws.Range("N48").ClearContents()

$ws.Range("H74").Value = 38592.6
$ws.Range("J74").Value = 38592.6
$ws.Range("L74").Value = 38592.6
$ws.Range("N74").Value = -40340.6

$ws.Range("H77").Value = 38592.6
$ws.Range("J77").Value = 38592.6
$ws.Range("L77").Value = 115777.8
$ws.Range("N77").Value = -124513.8

$ws.Range("H122").Value = 724.75
$ws.Range("I122").Value = 226
$ws.Range("J122").Value = 1223.5
$ws.Range("K122").Value = 678
$ws.Range("L122").Value = 3670.5
$ws.Range("M122").Value = 1772
$ws.Range("N122").Value = -8570.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 23842190
$ws.Range("J4").Value = 2554.4546
$ws.Range("L4").Value = 7663.3638
$ws.Range("N4").Value = -7887.3638

$ws.Range("H9").Value = 5000
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()

$ws.Range("H55").Value = 45620
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 57000
$ws.Range("K55").Value = 300
$ws.Range("L55").Value = 171000
$ws.Range("M55").Value = -123
$ws.Range("N55").Value = -171354

$ws.Range("H115").Value = 394
$ws.Range("I115").Value = 394
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 1182
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -7
$ws.Range("N115").ClearContents()

$ws.Range("H122").Value = 348.25
$ws.Range("J122").Value = 397.5
$ws.Range("L122").Value = 3577.5
$ws.Range("N122").Value = -8477.5

$ws.Range("H132").Value = 2129.1853
$ws.Range("I132").Value = 1716.8695
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 15451.8255
$ws.Range("L132").Value = 40500
$ws.Range("M132").Value = -12921.8255
$ws.Range("N132").Value = -45560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 24999
$ws.Range("I58").Value = 24999
$ws.Range("J58").Value = 24999
$ws.Range("K58").Value = 24999
$ws.Range("L58").Value = 24999
$ws.Range("M58").Value = -24722
$ws.Range("N58").Value = -25553

$ws.Range("H75").Value = 37000
$ws.Range("J75").Value = 37000
$ws.Range("L75").Value = 37000
$ws.Range("N75").Value = -38748

$ws.Range("H78").Value = 37000
$ws.Range("J78").Value = 37000
$ws.Range("L78").Value = 111000
$ws.Range("N78").Value = -119736

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1214.3334
$ws.Range("I22").Value = 1277.2
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 1277.2
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -982.2
$ws.Range("N22").Value = -1490

$ws.Range("H27").Value = 1214.3334
$ws.Range("I27").Value = 1277.2
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 1277.2
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = -1170.2
$ws.Range("N27").Value = -1114

$ws.Range("H122").Value = 2066.6667
$ws.Range("I122").Value = 2066.6667
$ws.Range("K122").Value = 6200.000100000001
$ws.Range("M122").Value = -3750.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H64").Value = 62249.75
$ws.Range("J64").Value = 62249.75
$ws.Range("L64").Value = 62249.75
$ws.Range("N64").Value = -62745.75

$ws.Range("H67").Value = 62249.75
$ws.Range("J67").Value = 62249.75
$ws.Range("L67").Value = 62249.75
$ws.Range("N67").Value = -63965.75

$ws.Range("H68").Value = 65000
$ws.Range("J68").Value = 65000
$ws.Range("L68").Value = 65000
$ws.Range("N68").Value = -66622

$ws.Range("H71").Value = 65000
$ws.Range("J71").Value = 65000
$ws.Range("L71").Value = 195000
$ws.Range("N71").Value = -203112

$ws.Range("H81").Value = 5137.5835
$ws.Range("J81").Value = 4243
$ws.Range("L81").Value = 8486
$ws.Range("N81").Value = -10608

$ws.Range("H84").Value = 5137.5835
$ws.Range("J84").Value = 4243
$ws.Range("L84").Value = 42430
$ws.Range("N84").Value = -53038

$ws.Range("H124").Value = 80000
$ws.Range("J124").Value = 80000
$ws.Range("L124").Value = 80000
$ws.Range("N124").Value = -89820

$ws.Range("H132").Value = 5137.875
$ws.Range("I132").Value = 5137.875
$ws.Range("K132").Value = 15413.625
$ws.Range("M132").Value = -12883.625
